$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(7,3,6,0),
    @(4,3,2,0),
    @(3,3,6,0),
    @(2,1,3,2),
    @(7,2,2,1),
    @(3,0,3,3),
    @(6,0,6,3),
    @(3,1,4,2),
    @(2,2,2,1),
    @(3,3,6,0),
    @(6,2,5,0),
    @(4,0,5,2),
    @(3,3,3,0),
    @(3,2,4,1),
    @(6,0,6,2),
    @(5,0,4,2),
    @(5,2,6,1),
    @(2,1,3,2),
    @(5,0,5,2),
    @(6,2,4,1),
    @(2,0,3,3),
    @(6,1,5,2),
    @(4,0,3,2),
    @(6,2,7,0),
    @(4,2,5,1)
)

$startRow = 1217
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
}

$excel.ActiveWindow.ScrollRow = 1214
$ws.Range("A1242").Select()
